$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.931.27"
$ws.Range("E2").Value = "  +6.54%  "

$ws.Range("D3").Value = "2.625.30"
$ws.Range("E3").Value = "  +8.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.14"
$ws.Range("E5").Value = "  +3.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.96"
$ws.Range("E6").Value = "  +1.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -3.53%  "

$ws.Range("D9").Value = "2.666.77"
$ws.Range("E9").Value = "  +9.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.46"
$ws.Range("E10").Value = "  +2.17%  "

$ws.Range("E11").Value = "  +5.17%  "

$ws.Range("E12").Value = "  +3.54%  "

$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").Value = "3.110.53"
$ws.Range("E14").Value = "  +9.92%  "

$ws.Range("D15").Value = "60.653.87"
$ws.Range("E15").Value = "  +6.14%  "

$ws.Range("E16").Value = "  +5.85%  "

$ws.Range("E17").Value = "  +6.27%  "

$ws.Range("D18").Value = "2.667.14"

$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.64"
$ws.Range("E20").Value = "  +6.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.54"
$ws.Range("E21").Value = "  +5.26%  "

$ws.Range("E22").Value = "  +4.79%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.38"
$ws.Range("E24").Value = "  +4.73%  "

$ws.Range("E25").Value = "  +5.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  +4.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").Value = "0.0₃0869"
$ws.Range("E28").Value = "  +11.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.63"
$ws.Range("E29").Value = "  +5.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.57"
$ws.Range("E31").Value = "  +5.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.95"
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("E33").Value = "  +3.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("E34").Value = "  +9.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.07"
$ws.Range("E35").Value = "  +7.65%  "

$ws.Range("E36").Value = "  +5.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "311.99"
$ws.Range("E37").Value = "  +12.04%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.860"
$ws.Range("E38").Value = "  +3.81%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  +9.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.840"
$ws.Range("E40").Value = "  +29.93%  "

$ws.Range("E41").Value = "  +7.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.35"
$ws.Range("E42").Value = "  +3.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.638"
$ws.Range("E43").Value = "  +6.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0577"
$ws.Range("E44").Value = "  +8.55%  "

$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.18"
$ws.Range("E46").Value = "  +14.97%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("E47").Value = "  -0.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.94"
$ws.Range("E48").Value = "  +9.30%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0237"
$ws.Range("E49").Value = "  +4.41%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.058.43"
$ws.Range("E50").Value = "  +8.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.27"
$ws.Range("E51").Value = "  +0.67%  "
